# Update the "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" worksheets, which contain duplicated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 158
    $ws.Range("F3").Value = 109
}
